$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 7.58
$ws.Range("A8").Value = -21.906
$ws.Range("A10").Value = -21.721
$ws.Range("A12").Value = -21.044
$ws.Range("B12").Value = 6.245
$ws.Range("D12").Value = -6.483
$ws.Range("D13").Value = -8.035999999999998
$ws.Range("B15").Value = 5.447
$ws.Range("B17").Value = 4.893000000000001
$ws.Range("A18").Value = -22.096
$ws.Range("D21").Value = -8.148
$ws.Range("D25").Value = -7.964
$ws.Range("B26").Value = 5.718999999999999
$ws.Range("B27").Value = 6.071
$ws.Range("B28").Value = 6.031999999999999
$ws.Range("D32").Value = -7.398999999999999
$ws.Range("D36").Value = -7.733
$ws.Range("A37").Value = -21.016
$ws.Range("B37").Value = 8.109
$ws.Range("D38").Value = -7.970000000000001
$ws.Range("D41").Value = -8.129000000000001
$ws.Range("B47").Value = 5.508
$ws.Range("D52").Value = -7.613000000000001
$ws.Range("A55").Value = -22.128
$ws.Range("D59").Value = -8.042999999999999
$ws.Range("B65").Value = 5.973999999999999
$ws.Range("D67").Value = -7.435
$ws.Range("A68").Value = -21.58300000000001
$ws.Range("B73").Value = 6.623
$ws.Range("A77").Value = -20.93
$ws.Range("A78").Value = -20.34200000000001
$ws.Range("A81").Value = -21.738
$ws.Range("A82").Value = -21.835
$ws.Range("B84").Value = 6.241000000000001
$ws.Range("D84").Value = -8.216000000000001
$ws.Range("B85").Value = 5.964
$ws.Range("D88").Value = -7.911000000000001
$ws.Range("D89").Value = -8.285
$ws.Range("B93").Value = 5.523000000000001
$ws.Range("B95").Value = 5.744
$ws.Range("D95").Value = -7.918000000000001
$ws.Range("B98").Value = 6.075
$ws.Range("B99").Value = 5.456999999999999
$ws.Range("B101").Value = 5.286
$ws.Range("D105").Value = -7.955999999999999
